# Update HealthProfessional mapping regarding the role.
# The "HealthProfessional.HealthProfessionalRole" mapping value that used
# to sit next to "EHDSHealthProfessional.role" in row 7 (column B) is
# relocated to a new row 11 (column B), leaving row 7 / column B empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the role mapping from its old location (B7) entirely, so the
# cell disappears from the sheet (not merely left blank).
$ws.Range("B7").Clear()

# Re-add it further down, in a new row 11, column B.
$ws.Range("B11").Value = "HealthProfessional.HealthProfessionalRole"

# Match the new active cell / selection recorded in the saved view.
$ws.Range("B11").Select()
